$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-04-29 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-04-30 Tuesday", 2)

# Update each answer cell in the single 20x5 table, in row-major order
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "8+32=40"
$t.Cell(1, 2).Range.Text = "86-85=1"
$t.Cell(1, 3).Range.Text = "19-17=2"
$t.Cell(1, 4).Range.Text = "51-33=18"
$t.Cell(1, 5).Range.Text = "98-69=29"

$t.Cell(2, 1).Range.Text = "99-84=15"
$t.Cell(2, 2).Range.Text = "87-6=81"
$t.Cell(2, 3).Range.Text = "5+71=76"
$t.Cell(2, 4).Range.Text = "31+27=58"
$t.Cell(2, 5).Range.Text = "9+58=67"

$t.Cell(3, 1).Range.Text = "45-13=32"
$t.Cell(3, 2).Range.Text = "17+36=53"
$t.Cell(3, 3).Range.Text = "99-81=18"
$t.Cell(3, 4).Range.Text = "65-53=12"
$t.Cell(3, 5).Range.Text = "98-97=1"

$t.Cell(4, 1).Range.Text = "78+8=86"
$t.Cell(4, 2).Range.Text = "35+19=54"
$t.Cell(4, 3).Range.Text = "50+0=50"
$t.Cell(4, 4).Range.Text = "34+11=45"
$t.Cell(4, 5).Range.Text = "52-26=26"

$t.Cell(5, 1).Range.Text = "77-44=33"
$t.Cell(5, 2).Range.Text = "21+58=79"
$t.Cell(5, 3).Range.Text = "61-13=48"
$t.Cell(5, 4).Range.Text = "80-14=66"
$t.Cell(5, 5).Range.Text = "95-52=43"

$t.Cell(6, 1).Range.Text = "3+66=69"
$t.Cell(6, 2).Range.Text = "69-30=39"
$t.Cell(6, 3).Range.Text = "0+78=78"
$t.Cell(6, 4).Range.Text = "96-77=19"
$t.Cell(6, 5).Range.Text = "14-8=6"

$t.Cell(7, 1).Range.Text = "19+19=38"
$t.Cell(7, 2).Range.Text = "98-24=74"
$t.Cell(7, 3).Range.Text = "56-44=12"
$t.Cell(7, 4).Range.Text = "42+8=50"
$t.Cell(7, 5).Range.Text = "31+2=33"

$t.Cell(8, 1).Range.Text = "14+36=50"
$t.Cell(8, 2).Range.Text = "75-16=59"
$t.Cell(8, 3).Range.Text = "50-22=28"
$t.Cell(8, 4).Range.Text = "70+3=73"
$t.Cell(8, 5).Range.Text = "96-3=93"

$t.Cell(9, 1).Range.Text = "37+16=53"
$t.Cell(9, 2).Range.Text = "39+45=84"
$t.Cell(9, 3).Range.Text = "89-44=45"
$t.Cell(9, 4).Range.Text = "77-76=1"
$t.Cell(9, 5).Range.Text = "61-42=19"

$t.Cell(10, 1).Range.Text = "56-0=56"
$t.Cell(10, 2).Range.Text = "6+14=20"
$t.Cell(10, 3).Range.Text = "5+8=13"
$t.Cell(10, 4).Range.Text = "69-53=16"
$t.Cell(10, 5).Range.Text = "95-85=10"

$t.Cell(11, 1).Range.Text = "72-62=10"
$t.Cell(11, 2).Range.Text = "91-34=57"
$t.Cell(11, 3).Range.Text = "37+58=95"
$t.Cell(11, 4).Range.Text = "72-37=35"
$t.Cell(11, 5).Range.Text = "86-23=63"

$t.Cell(12, 1).Range.Text = "87-58=29"
$t.Cell(12, 2).Range.Text = "67-40=27"
$t.Cell(12, 3).Range.Text = "73+26=99"
$t.Cell(12, 4).Range.Text = "53-20=33"
$t.Cell(12, 5).Range.Text = "24+70=94"

$t.Cell(13, 1).Range.Text = "87-83=4"
$t.Cell(13, 2).Range.Text = "93-52=41"
$t.Cell(13, 3).Range.Text = "23+34=57"
$t.Cell(13, 4).Range.Text = "87-20=67"
$t.Cell(13, 5).Range.Text = "92-35=57"

$t.Cell(14, 1).Range.Text = "59-8=51"
$t.Cell(14, 2).Range.Text = "9+49=58"
$t.Cell(14, 3).Range.Text = "21+37=58"
$t.Cell(14, 4).Range.Text = "79-19=60"
$t.Cell(14, 5).Range.Text = "1+41=42"

$t.Cell(15, 1).Range.Text = "34+16=50"
$t.Cell(15, 2).Range.Text = "66+6=72"
$t.Cell(15, 3).Range.Text = "64-50=14"
$t.Cell(15, 4).Range.Text = "25+20=45"
$t.Cell(15, 5).Range.Text = "49+23=72"

$t.Cell(16, 1).Range.Text = "59-53=6"
$t.Cell(16, 2).Range.Text = "87-15=72"
$t.Cell(16, 3).Range.Text = "83-62=21"
$t.Cell(16, 4).Range.Text = "54-47=7"
$t.Cell(16, 5).Range.Text = "76-38=38"

$t.Cell(17, 1).Range.Text = "26-6=20"
$t.Cell(17, 2).Range.Text = "65-57=8"
$t.Cell(17, 3).Range.Text = "23-6=17"
$t.Cell(17, 4).Range.Text = "79+18=97"
$t.Cell(17, 5).Range.Text = "69+1=70"

$t.Cell(18, 1).Range.Text = "91-90=1"
$t.Cell(18, 2).Range.Text = "6-2=4"
$t.Cell(18, 3).Range.Text = "60-0=60"
$t.Cell(18, 4).Range.Text = "15-4=11"
$t.Cell(18, 5).Range.Text = "67-11=56"

$t.Cell(19, 1).Range.Text = "39+28=67"
$t.Cell(19, 2).Range.Text = "10+21=31"
$t.Cell(19, 3).Range.Text = "47+27=74"
$t.Cell(19, 4).Range.Text = "44-44=0"
$t.Cell(19, 5).Range.Text = "45+6=51"

$t.Cell(20, 1).Range.Text = "68+1=69"
$t.Cell(20, 2).Range.Text = "47+37=84"
$t.Cell(20, 3).Range.Text = "85-28=57"
$t.Cell(20, 4).Range.Text = "71-68=3"
$t.Cell(20, 5).Range.Text = "70-33=37"

